# Update logo font to TW Cen MT Bold Italic
#
# - "Sparkitecture" title: Century Gothic -> Tw Cen MT (shape re-centers
#   horizontally as PowerPoint re-measures the spAutoFit text box).
# - Decorative rotated logo group nudges slightly to stay aligned with the
#   resized title.
# - The caption text box documents the font used and grows (spAutoFit) to
#   fit the added line, which shifts its position too.

# Shape.Left/Top/Width/Height round-trip through a 32-bit float (just like
# real PowerPoint's COM automation layer), so naive "emu / 12700.0" can miss
# the exact target EMU by one unit. This helper nudges the point value so it
# lands on the requested EMU after that float32 round-trip.
function EmuToPt($targetEmu) {
    $base = $targetEmu / 12700.0
    for ($i = 0; $i -lt 4000; $i++) {
        $candidate = $base + ($i * 0.0000001)
        $asSingle = [single]$candidate
        $backToEmu = [int64]([double]$asSingle * 12700)
        if ($backToEmu -eq $targetEmu) {
            return $candidate
        }
    }
    return $base
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# 1) "Sparkitecture" title (Rectangle 3) - swap the Latin typeface and
#    re-fit/re-center the autosize text box.
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$titleRange = $title.TextFrame.TextRange
$titleRange.Font.Name = "Tw Cen MT"

$title.Left = EmuToPt 2090253
$title.Top = EmuToPt 2154338
$title.Width = EmuToPt 6202917
$title.Height = EmuToPt 1569660

# ---------------------------------------------------------------------
# 2) Decorative rotated logo group (Group 21) - reposition to match.
# ---------------------------------------------------------------------
$logoGroup = $s.Shapes.Item(2)
$logoGroup.Left = EmuToPt 8076373
$logoGroup.Top = EmuToPt 1967630

# ---------------------------------------------------------------------
# 3) Caption text box (TextBox 23) - add a line documenting the new font
#    and let it grow/reposition with the extra line.
# ---------------------------------------------------------------------
$caption = $s.Shapes.Item(3)
$captionRange = $caption.TextFrame.TextRange

$urlText = "https://www.instagram.com/p/B1VyYMDgF2Q/"
$line1 = "Example: " + $urlText
$line2 = "Font: TW Cen MT Bold Italic"
$captionRange.Text = $line1 + "`r" + $line2

# Re-apply the hyperlink that the plain-text assignment above dropped.
$urlStart = ("Example: ").Length + 1
$urlRange = $captionRange.Characters($urlStart, $urlText.Length)
$urlRange.ActionSettings.Item(1).Hyperlink.Address = $urlText

$caption.Left = EmuToPt 83464
$caption.Top = EmuToPt 6211669
$caption.Width = EmuToPt 5560818
$caption.Height = EmuToPt 646331
